$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50; existing rows 50-82 shift down to 51-83.
$ws.Rows("50:50").Insert()

# Populate the newly inserted row 50 with the new weekly price record.
$ws.Cells.Item(50,1).Value  = 8
$ws.Cells.Item(50,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(50,3).Value  = "Coquimbo"
$ws.Cells.Item(50,4).Value  = 44673
$ws.Cells.Item(50,5).Value  = 4
$ws.Cells.Item(50,6).Value  = 100112052
$ws.Cells.Item(50,7).Value  = "Albahaca"
$ws.Cells.Item(50,8).Value  = "Sin especificar"
$ws.Cells.Item(50,9).Value  = "Primera"
$ws.Cells.Item(50,10).Value = 800
$ws.Cells.Item(50,11).Value = 5000
$ws.Cells.Item(50,12).Value = 5500
$ws.Cells.Item(50,13).Value = 5250
$ws.Cells.Item(50,14).Value = "`$/docena de matas"
$ws.Cells.Item(50,15).Value = "Provincia del Elquí"
$ws.Cells.Item(50,16).Value = 875
$ws.Cells.Item(50,17).Value = 6
$ws.Cells.Item(50,18).Value = "Hortaliza"

# Make sure the date column keeps the same number-format style used by the
# other rows (style index matching the "D" column date cells).
$ws.Cells.Item(50,4).NumberFormat = $ws.Cells.Item(51,4).NumberFormat
